$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H40").Value = 1266.6666
$ws_ALC.Range("I40").Value = 1226.9231
$ws_ALC.Range("J40").Value = 1414.2858
$ws_ALC.Range("K40").Value = 1226.9231
$ws_ALC.Range("L40").Value = 1414.2858
$ws_ALC.Range("M40").Value = -1051.9231
$ws_ALC.Range("N40").Value = -1764.2858
$ws_ALC.Range("H87").Value = 12177.944
$ws_ALC.Range("J87").Value = 13165.5
$ws_ALC.Range("L87").Value = 13165.5
$ws_ALC.Range("N87").Value = -15661.5
$ws_ALC.Range("H90").Value = 12177.944
$ws_ALC.Range("J90").Value = 13165.5
$ws_ALC.Range("L90").Value = 39496.5
$ws_ALC.Range("N90").Value = -51976.5
$ws_ALC.Range("H97").Value = 1065.5
$ws_ALC.Range("J97").Value = 1881
$ws_ALC.Range("L97").Value = 5643
$ws_ALC.Range("N97").Value = -6635
$ws_ALC.Range("H99").Value = 700.7778
$ws_ALC.Range("I99").Value = 254.66667
$ws_ALC.Range("J99").Value = 1593
$ws_ALC.Range("K99").Value = 764.00001
$ws_ALC.Range("L99").Value = 4779
$ws_ALC.Range("M99").Value = 733.99999
$ws_ALC.Range("N99").Value = -7775
$ws_ALC.Range("H115").Value = 1685
$ws_ALC.Range("I115").Value = 299.16666
$ws_ALC.Range("K115").Value = 897.4999799999999
$ws_ALC.Range("M115").Value = 669.5000200000001

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H74").Value = 1915.15
$ws_ARM.Range("I74").Value = 977.1539
$ws_ARM.Range("J74").Value = 3657.1428
$ws_ARM.Range("K74").Value = 977.1539
$ws_ARM.Range("L74").Value = 3657.1428
$ws_ARM.Range("M74").Value = -103.1539
$ws_ARM.Range("N74").Value = -5405.1428
$ws_ARM.Range("H77").Value = 1915.15
$ws_ARM.Range("I77").Value = 977.1539
$ws_ARM.Range("J77").Value = 3657.1428
$ws_ARM.Range("K77").Value = 4885.7695
$ws_ARM.Range("L77").Value = 18285.714
$ws_ARM.Range("M77").Value = -517.7695000000003
$ws_ARM.Range("N77").Value = -27021.714
$ws_ARM.Range("H122").Value = 2394.2307
$ws_ARM.Range("I122").Value = 2115.0527
$ws_ARM.Range("J122").Value = 3152
$ws_ARM.Range("K122").Value = 6345.158100000001
$ws_ARM.Range("L122").Value = 9456
$ws_ARM.Range("M122").Value = -3895.158100000001
$ws_ARM.Range("N122").Value = -14356
$ws_ARM.Range("H130").Value = 23095.334
$ws_ARM.Range("J130").Value = 23095.334
$ws_ARM.Range("L130").Value = 23095.334
$ws_ARM.Range("N130").Value = -33135.334
$ws_ARM.Range("H132").Value = 2532.7869
$ws_ARM.Range("I132").Value = 2248.449
$ws_ARM.Range("K132").Value = 6745.347
$ws_ARM.Range("M132").Value = -4215.347

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H134").Value = 1450.9231
$ws_BSM.Range("I134").Value = 1290.2858
$ws_BSM.Range("J134").Value = 2125.6
$ws_BSM.Range("K134").Value = 3870.8574
$ws_BSM.Range("L134").Value = 6376.799999999999
$ws_BSM.Range("M134").Value = -1335.8574
$ws_BSM.Range("N134").Value = -11446.8

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H31").Value = 5421
$ws_CRP.Range("I31").Value = 6182.4
$ws_CRP.Range("J31").Value = 4835.3076
$ws_CRP.Range("K31").Value = 6182.4
$ws_CRP.Range("L31").Value = 4835.3076
$ws_CRP.Range("M31").Value = -5887.4
$ws_CRP.Range("N31").Value = -5425.3076
$ws_CRP.Range("H34").Value = 5421
$ws_CRP.Range("I34").Value = 6182.4
$ws_CRP.Range("J34").Value = 4835.3076
$ws_CRP.Range("K34").Value = 6182.4
$ws_CRP.Range("L34").Value = 4835.3076
$ws_CRP.Range("M34").Value = -5980.4
$ws_CRP.Range("N34").Value = -5239.3076

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H5").Value = 23094.305
$ws_CUL.Range("I5").Value = 384.16
$ws_CUL.Range("J5").Value = 50130.19
$ws_CUL.Range("K5").Value = 1152.48
$ws_CUL.Range("L5").Value = 150390.57
$ws_CUL.Range("M5").Value = -1040.48
$ws_CUL.Range("N5").Value = -150614.57
$ws_CUL.Range("H55").Value = 2537.75
$ws_CUL.Range("J55").Value = 2761.5386
$ws_CUL.Range("L55").Value = 8284.6158
$ws_CUL.Range("N55").Value = -8638.6158
$ws_CUL.Range("H87").Value = 7500
$ws_CUL.Range("I87").Value = 1666.6666
$ws_CUL.Range("J87").Value = 25000
$ws_CUL.Range("K87").Value = 4999.9998
$ws_CUL.Range("L87").Value = 75000
$ws_CUL.Range("M87").Value = -3751.9998
$ws_CUL.Range("N87").Value = -77496
$ws_CUL.Range("H90").Value = 7500
$ws_CUL.Range("I90").Value = 1666.6666
$ws_CUL.Range("J90").Value = 25000
$ws_CUL.Range("K90").Value = 14999.9994
$ws_CUL.Range("L90").Value = 225000
$ws_CUL.Range("M90").Value = -8759.999400000001
$ws_CUL.Range("N90").Value = -237480
$ws_CUL.Range("H103").Value = 2738.0476
$ws_CUL.Range("J103").Value = 4002.0715
$ws_CUL.Range("L103").Value = 12006.2145
$ws_CUL.Range("N103").Value = -13764.2145
$ws_CUL.Range("H124").Value = 5553.75
$ws_CUL.Range("I124").Value = 2915
$ws_CUL.Range("J124").Value = 6433.3335
$ws_CUL.Range("K124").Value = 8745
$ws_CUL.Range("L124").Value = 19300.0005
$ws_CUL.Range("M124").Value = -3835
$ws_CUL.Range("N124").Value = -29120.0005
$ws_CUL.Range("H131").Value = 897.7931
$ws_CUL.Range("I131").Value = 544.4545000000001
$ws_CUL.Range("J131").Value = 980.4894
$ws_CUL.Range("K131").Value = 1633.3635
$ws_CUL.Range("L131").Value = 2941.4682
$ws_CUL.Range("M131").Value = 3406.6365
$ws_CUL.Range("N131").Value = -13021.4682
$ws_CUL.Range("H135").Value = 23094.305
$ws_CUL.Range("I135").Value = 384.16
$ws_CUL.Range("J135").Value = 50130.19
$ws_CUL.Range("K135").Value = 3457.44
$ws_CUL.Range("L135").Value = 451171.71
$ws_CUL.Range("M135").Value = -922.4400000000001
$ws_CUL.Range("N135").Value = -456241.71

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H80").Value = 2999.5
$ws_GSM.Range("I80").Value = 2916.6667
$ws_GSM.Range("J80").Value = 3248
$ws_GSM.Range("K80").Value = 2916.6667
$ws_GSM.Range("L80").Value = 3248
$ws_GSM.Range("M80").Value = -1918.6667
$ws_GSM.Range("N80").Value = -5244
$ws_GSM.Range("H83").Value = 2999.5
$ws_GSM.Range("I83").Value = 2916.6667
$ws_GSM.Range("J83").Value = 3248
$ws_GSM.Range("K83").Value = 14583.3335
$ws_GSM.Range("L83").Value = 16240
$ws_GSM.Range("M83").Value = -9591.333500000001
$ws_GSM.Range("N83").Value = -26224
$ws_GSM.Range("H132").Value = 1808.7858
$ws_GSM.Range("I132").Value = 1368.8928
$ws_GSM.Range("J132").Value = 2688.5715
$ws_GSM.Range("K132").Value = 4106.678400000001
$ws_GSM.Range("L132").Value = 8065.7145
$ws_GSM.Range("M132").Value = -1576.678400000001
$ws_GSM.Range("N132").Value = -13125.7145

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H40").Value = 49500
$ws_WVR.Range("I40").Value = 0
$ws_WVR.Range("J40").Value = 49500
$ws_WVR.Range("K40").Value = 0
$ws_WVR.Range("L40").Value = 49500
$ws_WVR.Range("M40").ClearContents()
$ws_WVR.Range("N40").Value = -49798
$ws_WVR.Range("H42").Value = 70000
$ws_WVR.Range("J42").Value = 70000
$ws_WVR.Range("L42").Value = 70000
$ws_WVR.Range("N42").Value = -70756
